$wb = $excel.ActiveWorkbook

# Sheet references (index is 1-based, matches workbook.xml sheet order)
$wsPrep = $wb.Worksheets.Item(1)   # "Prep phase"
$wsAI   = $wb.Worksheets.Item(3)   # "AI phase "
$wsPost = $wb.Worksheets.Item(4)   # "Post training phase"

# --- Prep phase: add row 2 content ---
$wsPrep.Range("A2").Value = "Images are to be manually annotated"
$wsPrep.Range("B2").Value = "Using labelimg tool manually annotate and fetch xmin,ymin of img."

# Widen columns A and B (closest reachable widths to the target 34.77734375 / 59.109375)
$wsPrep.Columns.Item(1).ColumnWidth = 34
$wsPrep.Columns.Item(2).ColumnWidth = 58.333333333333336

# Selection ends on C2 for this sheet, and it will not be the active tab
$wsPrep.Range("C2").Select()

# --- AI phase: add row 2 content ---
$wsAI.Range("A2").Value = "training the model"
$wsAI.Range("B2").Value = "thrain the ml model with annoted images"

# Selection ends on B2 for this sheet
$wsAI.Range("B2").Select()

# --- Post training phase: add row 2 content ---
$wsPost.Range("A2").Value = "create inference graph"
$wsPost.Range("B2").Value = "save the trained model and use it to predict o/p"

# Widen column A (closest reachable width to the target 20.21875)
$wsPost.Columns.Item(1).ColumnWidth = 19.333333333333332

# Selection ends on B12 for this sheet
$wsPost.Range("B12").Select()

# Make "Post training phase" the active sheet/tab (activeTab=3, tabSelected moves here)
$wsPost.Activate()
